$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.591.82'
$ws.Range("E2").Value = '  +3.53%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.072.58'
$ws.Range("E3").Value = '  +2.74%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '550.68'
$ws.Range("E5").Value = '  +2.65%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.02'
$ws.Range("E6").Value = '  +7.59%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.068.39'
$ws.Range("E8").Value = '  +2.78%  '

$ws.Range("E9").Value = '  +1.55%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.50'
$ws.Range("E10").Value = '  +6.17%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.152'
$ws.Range("E11").Value = '  +3.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.458'
$ws.Range("E12").Value = '  +3.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000228'
$ws.Range("E13").Value = '  +3.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.93'
$ws.Range("E14").Value = '  +3.85%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.569.78'
$ws.Range("E15").Value = '  +2.53%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.537.10'
$ws.Range("E16").Value = '  +3.26%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.073.19'
$ws.Range("E17").Value = '  +2.31%  '

$ws.Range("E18").Value = '  -0.72%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.79'
$ws.Range("E19").Value = '  +3.05%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '484.12'
$ws.Range("E20").Value = '  +3.90%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.91'
$ws.Range("E21").Value = '  +5.67%  '

$ws.Range("E22").Value = '  +1.34%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.32'
$ws.Range("E23").Value = '  +6.33%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.11'
$ws.Range("E24").Value = '  +0.79%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.85'
$ws.Range("E25").Value = '  +8.00%  '

$ws.Range("E26").Value = '  +0.15%  '

$ws.Range("E27").Value = '  +4.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.95'
$ws.Range("E28").Value = '  +3.48%  '

$ws.Range("E29").Value = '  +8.36%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.22'
$ws.Range("E31").Value = '  +2.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.18'
$ws.Range("E32").Value = '  +2.02%  '

$ws.Range("E33").Value = '  +8.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.71'
$ws.Range("E34").Value = '  +5.20%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '55.59'
$ws.Range("E35").Value = '  +0.62%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.02'
$ws.Range("E36").Value = '  +2.75%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '468.75'
$ws.Range("E37").Value = '  +4.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0823'
$ws.Range("E38").Value = '  +5.25%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0399'
$ws.Range("E39").Value = '  +4.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.030.09'
$ws.Range("E40").Value = '  -4.09%  '

$ws.Range("E41").Value = '  +0.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.23'
$ws.Range("E42").Value = '  +2.02%  '

$ws.Range("E43").Value = '  +6.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '27.74'
$ws.Range("E44").Value = '  +5.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.256'
$ws.Range("E45").Value = '  +5.87%  '

$ws.Range("E47").Value = '  +3.36%  '

$ws.Range("E48").Value = '  +2.88%  '

$ws.Range("E51").Value = '  +4.61%  '

# Row 49/50: PEPE and Monero swap positions with updated price/volume data
$sub3 = [char]0x2083

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '116.75'
$ws.Range("E49").Value = '  -1.15%  '

$ws.Range("B50").Value = 'PEPE'
$ws.Range("C50").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D50").Value = "0.0$sub3" + "0511"
$ws.Range("E50").Value = '  +3.86%  '
